$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark (currently at the end of the
#    "This question was very easy to implement." paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Append a new run of text after the existing hyperlink in the
#    References paragraph, styled like a hyperlink.
$hyperlink = $d.Hyperlinks(1)
$tail = $hyperlink.Range
$tail.Collapse(0)  # wdCollapseEnd
$tail.InsertAfter(" http://www.gnu.org/software/bison/manual/bison.html")
$tail.Style = "Hyperlink"

# 3. Re-insert the "_GoBack" bookmark at the very end of the document.
$end = $d.Content
$end.Collapse(0)  # wdCollapseEnd
$d.Bookmarks.Add("_GoBack", $end)

$word.ActiveDocument.Save()
